# Apply the edits described by the diff to All_missions.xlsx (Tabelle1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# 1. Update the big description box (merged M4:W13) for the "Main Tasks" table:
#    it now shows the "Mission object detection" description text instead of the
#    generic board-layout hint text.
$ws.Range("M4").Value = "Erkenne während der Linienverfolgung, ob sich Objekte sichtbar oder fühlbar vor Zumi auf der Straße befinden (schwarze Linie). Gib ein Geräusch ab, wenn du ein Objekt erkennst, und halte den Roboter an, wenn er nahe genug an dem Objekt ist. Sobald Zumi angehalten hat, nimmst du das Objekt mit der Hand von der Straße und fährst weiter. Zählen Sie die Anzahl der erkannten Objekte. Wir unterschei-den nicht zwischen den Objekttypen, sondern konzentrieren uns nur darauf, ob ein Objekt erkannt wird. Alle Objekte sind groß genug, um vom vorderen IR-Sensor erkannt zu werden."

# 2. Row 7 ("Mission object detection" status row in first table): status moves
#    from "Semi" (I7) to "Done" (H7).
$ws.Range("H7").Value = 1
$ws.Range("I7").ClearContents()

# 3. Row 30 (Surprise task row in second table): status moves from "Unfin." (J30)
#    to "Done" (H30).
$ws.Range("H30").Value = 1
$ws.Range("J30").ClearContents()

# Force a full recalculation so the SUM()/percentage formulas in rows 16-17,
# 39-40 and 62 (and K62) pick up the new status values.
$excel.CalculateFullRebuild()

# 4. Update the sheet view / current selection to match where the author was
#    working when the file was saved.
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H8").Select()
